$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.132844
$ws.Range("H2").Value = 3.398531999999999
$ws.Range("I2").Value = 0.6473844372037293
$ws.Range("J2").Value = 0.6473844372037293
$ws.Range("M2").Value = 30.52246933333333
$ws.Range("N2").Value = 91.567408
$ws.Range("O2").Value = 0.1058764512547768
$ws.Range("P2").Value = 0.1058764512547769
$ws.Range("Q2").Value = 34.57719624945066
$ws.Range("R2").Value = 311.194766245056
$ws.Range("S2").Value = 0.06854276680870179
$ws.Range("T2").Value = 0.06854276680870179

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.132844
$ws.Range("H3").Value = 3.398531999999999
$ws.Range("I3").Value = 0.6473844372037293
$ws.Range("J3").Value = 0.6473844372037293
$ws.Range("O3").Value = 0.001067503492562006
$ws.Range("P3").Value = 0.001067503492562006
$ws.Range("Q3").Value = 0.3486259439359999
$ws.Range("R3").Value = 3.137633495423999
$ws.Range("S3").Value = 0.0006910851477452694
$ws.Range("T3").Value = 0.0006910851477452695

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.132844
$ws.Range("H4").Value = 3.398531999999999
$ws.Range("I4").Value = 0.6473844372037293
$ws.Range("J4").Value = 0.6473844372037293
$ws.Range("M4").Value = 47.57542166666667
$ws.Range("N4").Value = 142.726265
$ws.Range("O4").Value = 0.1650297935598315
$ws.Range("P4").Value = 0.1650297935598315
$ws.Range("Q4").Value = 53.89553098255332
$ws.Range("R4").Value = 485.0597788429799
$ws.Range("S4").Value = 0.1068377200255792
$ws.Range("T4").Value = 0.1068377200255792

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.132844
$ws.Range("H5").Value = 3.398531999999999
$ws.Range("I5").Value = 0.6473844372037293
$ws.Range("J5").Value = 0.6473844372037293
$ws.Range("M5").Value = 209.8781993333333
$ws.Range("N5").Value = 629.634598
$ws.Range("O5").Value = 0.7280262516928295
$ws.Range("P5").Value = 0.7280262516928296
$ws.Range("Q5").Value = 237.7592588455706
$ws.Range("R5").Value = 2139.833329610135
$ws.Range("S5").Value = 0.471312865221703
$ws.Range("T5").Value = 0.4713128652217031

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6170343333333334
$ws.Range("H6").Value = 1.851103
$ws.Range("I6").Value = 0.3526155627962707
$ws.Range("J6").Value = 0.3526155627962707
$ws.Range("M6").Value = 30.52246933333333
$ws.Range("N6").Value = 91.567408
$ws.Range("O6").Value = 0.1058764512547768
$ws.Range("P6").Value = 0.1058764512547769
$ws.Range("Q6").Value = 18.83341151678044
$ws.Range("R6").Value = 169.500703651024
$ws.Range("S6").Value = 0.03733368444607506
$ws.Range("T6").Value = 0.03733368444607506

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.6170343333333334
$ws.Range("H7").Value = 1.851103
$ws.Range("I7").Value = 0.3526155627962707
$ws.Range("J7").Value = 0.3526155627962707
$ws.Range("O7").Value = 0.001067503492562006
$ws.Range("P7").Value = 0.001067503492562006
$ws.Range("Q7").Value = 0.1898886138773333
$ws.Range("R7").Value = 1.708997524896
$ws.Range("S7").Value = 0.0003764183448167361
$ws.Range("T7").Value = 0.0003764183448167362

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6170343333333334
$ws.Range("H8").Value = 1.851103
$ws.Range("I8").Value = 0.3526155627962707
$ws.Range("J8").Value = 0.3526155627962707
$ws.Range("M8").Value = 47.57542166666667
$ws.Range("N8").Value = 142.726265
$ws.Range("O8").Value = 0.1650297935598315
$ws.Range("P8").Value = 0.1650297935598315
$ws.Range("Q8").Value = 29.35566859114389
$ws.Range("R8").Value = 264.201017320295
$ws.Range("S8").Value = 0.05819207353425235
$ws.Range("T8").Value = 0.05819207353425235

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6170343333333334
$ws.Range("H9").Value = 1.851103
$ws.Range("I9").Value = 0.3526155627962707
$ws.Range("J9").Value = 0.3526155627962707
$ws.Range("M9").Value = 209.8781993333333
$ws.Range("N9").Value = 629.634598
$ws.Range("O9").Value = 0.7280262516928295
$ws.Range("P9").Value = 0.7280262516928296
$ws.Range("Q9").Value = 129.5020548068438
$ws.Range("R9").Value = 1165.518493261594
$ws.Range("S9").Value = 0.2567133864711265
$ws.Range("T9").Value = 0.2567133864711265
